$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("S3").Value = 2023
$ws.Range("S4").Value = 130.1
$ws.Range("S5").Value = 123
$ws.Range("S6").Value = 1488
$ws.Range("S7").Value = 931
$ws.Range("S8").Value = 1179.1
$ws.Range("S9").Value = 56.2
$ws.Range("S10").Value = 13.4
$ws.Range("S11").Value = 66.8
$ws.Range("S12").Value = 6.1
$ws.Range("S13").Value = 64.2
$ws.Range("S14").Value = 8.8
